# Views.xlsx edit: add "all/" (PublicCookieCreationsView) data row to the
# api/v1/cookiecreations/ section, and add a brand-new api/v1/mycookies/
# section (AllUsersCreations / AllUsersFavorites / AddCreationToFavorites).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- api/v1/cookiecreations/ section (row 15) gains a data row ---------
# (A15 already holds "api/v1/cookiecreations/" - leave it alone)
$ws.Range("B11").Value = "all/"
$ws.Range("C15").Value = "PublicCookieCreationsView"
$ws.Range("E15").Value = "json: all previously bought cookie creations (where prev_purchased = t)"

# --- new api/v1/mycookies/ section (rows 11-13) -------------------------
$ws.Range("B13").Value = "add/"
$ws.Range("B12").Value = "favorites/"
$ws.Range("A11").Value = "api/v1/mycookies/"
$ws.Range("C11").Value = "AllUsersCreations"
$ws.Range("C12").Value = "AllUsersFavorites"
$ws.Range("C13").Value = "AddCreationToFavorites"
$ws.Range("E11").Value = "json: all personal creations"
$ws.Range("E12").Value = "json: all favorites"

# --- remaining (reused / pre-existing) string cells ---------------------
$ws.Range("D11").Value = "get`n--auth header`n--json: empty"
$ws.Range("D12").Value = "get`n--auth header`n--json: empty"
$ws.Range("D15").Value = "get`n--auth header"
$ws.Range("F11").Value = "200: ok"
$ws.Range("F12").Value = "200: ok"
$ws.Range("F15").Value = "200: ok"
$ws.Range("B15").Value = "all/"

# --- wrap text (reuses the existing wrap-text cell style) ---------------
$ws.Range("D11").WrapText = $true
$ws.Range("D12").WrapText = $true
$ws.Range("D15").WrapText = $true
$ws.Range("E15").WrapText = $true

# --- row heights (match the sheet's existing "16px per wrapped line") ---
$ws.Rows.Item(11).RowHeight = 48
$ws.Rows.Item(12).RowHeight = 48
$ws.Rows.Item(15).RowHeight = 32

# --- column widths (C widens for "AddCreationToFavorites", E widens for
#     the long "previously bought cookie creations" response text) ------
$ws.Columns.Item(3).ColumnWidth = 22.41796875
$ws.Columns.Item(5).ColumnWidth = 35.25390625

# --- selection moves to F13 ---------------------------------------------
$ws.Range("F13").Select()
